{"js": "// Replace each multiplication expression in the practice table with the\n// values from the updated worksheet (20 rows x 5 columns, row-major order,\n// same order as the cells already appear in the document).\nconst newValues = [\n  [\"35\u00d749=\", \"57\u00d723=\", \"72\u00d743=\", \"49\u00d799=\", \"94\u00d723=\"],\n  [\"42\u00d715=\", \"52\u00d795=\", \"53\u00d771=\", \"54\u00d755=\", \"38\u00d734=\"],\n  [\"92\u00d772=\", \"16\u00d776=\", \"87\u00d753=\", \"41\u00d733=\", \"85\u00d748=\"],\n  [\"72\u00d775=\", \"37\u00d712=\", \"30\u00d792=\", \"15\u00d788=\", \"74\u00d744=\"],\n  [\"78\u00d723=\", \"49\u00d723=\", \"86\u00d789=\", \"18\u00d764=\", \"78\u00d716=\"],\n  [\"40\u00d751=\", \"52\u00d7100=\", \"44\u00d774=\", \"61\u00d768=\", \"71\u00d782=\"],\n  [\"80\u00d775=\", \"17\u00d752=\", \"85\u00d756=\", \"18\u00d771=\", \"40\u00d742=\"],\n  [\"11\u00d716=\", \"16\u00d720=\", \"14\u00d719=\", \"100\u00d749=\", \"28\u00d756=\"],\n  [\"20\u00d754=\", \"57\u00d730=\", \"43\u00d748=\", \"14\u00d790=\", \"41\u00d760=\"],\n  [\"20\u00d717=\", \"30\u00d743=\", \"39\u00d776=\", \"55\u00d712=\", \"64\u00d751=\"],\n  [\"25\u00d743=\", \"38\u00d753=\", \"25\u00d754=\", \"47\u00d712=\", \"48\u00d772=\"],\n  [\"64\u00d756=\", \"47\u00d791=\", \"22\u00d729=\", \"13\u00d780=\", \"60\u00d793=\"],\n  [\"44\u00d722=\", \"25\u00d752=\", \"34\u00d730=\", \"35\u00d753=\", \"80\u00d784=\"],\n  [\"77\u00d752=\", \"49\u00d754=\", \"37\u00d762=\", \"75\u00d758=\", \"14\u00d771=\"],\n  [\"84\u00d728=\", \"20\u00d716=\", \"73\u00d745=\", \"74\u00d788=\", \"34\u00d711=\"],\n  [\"96\u00d745=\", \"36\u00d775=\", \"28\u00d736=\", \"40\u00d749=\", \"57\u00d793=\"],\n  [\"62\u00d718=\", \"44\u00d743=\", \"70\u00d785=\", \"73\u00d793=\", \"89\u00d756=\"],\n  [\"76\u00d721=\", \"31\u00d772=\", \"14\u00d763=\", \"33\u00d717=\", \"41\u00d710=\"],\n  [\"68\u00d771=\", \"74\u00d789=\", \"10\u00d733=\", \"74\u00d784=\", \"94\u00d736=\"],\n  [\"10\u00d782=\", \"18\u00d767=\", \"43\u00d754=\", \"45\u00d723=\", \"77\u00d782=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst actualColumnCount = table.values.length > 0 ? table.values[0].length : 0;\nif (table.rowCount !== newValues.length || actualColumnCount !== newValues[0].length) {\n  throw new Error(\n    `Table shape ${table.rowCount}x${actualColumnCount} does not match expected ` +\n    `${newValues.length}x${newValues[0].length}.`\n  );\n}\n\n// Word.Table.values lets us overwrite every cell's text in one shot while\n// leaving paragraph/run formatting (font, size, alignment, etc.) untouched.\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Each line below is one table row (5 comma-free values joined by a\n# delimiter) holding the new \"A x B=\" text for that row's 5 cells, in the\n# same row-major order the cells already appear in the document's single\n# table (20 rows x 5 columns).\n$newRows = @\"\n35\u00d749=|57\u00d723=|72\u00d743=|49\u00d799=|94\u00d723=\n42\u00d715=|52\u00d795=|53\u00d771=|54\u00d755=|38\u00d734=\n92\u00d772=|16\u00d776=|87\u00d753=|41\u00d733=|85\u00d748=\n72\u00d775=|37\u00d712=|30\u00d792=|15\u00d788=|74\u00d744=\n78\u00d723=|49\u00d723=|86\u00d789=|18\u00d764=|78\u00d716=\n40\u00d751=|52\u00d7100=|44\u00d774=|61\u00d768=|71\u00d782=\n80\u00d775=|17\u00d752=|85\u00d756=|18\u00d771=|40\u00d742=\n11\u00d716=|16\u00d720=|14\u00d719=|100\u00d749=|28\u00d756=\n20\u00d754=|57\u00d730=|43\u00d748=|14\u00d790=|41\u00d760=\n20\u00d717=|30\u00d743=|39\u00d776=|55\u00d712=|64\u00d751=\n25\u00d743=|38\u00d753=|25\u00d754=|47\u00d712=|48\u00d772=\n64\u00d756=|47\u00d791=|22\u00d729=|13\u00d780=|60\u00d793=\n44\u00d722=|25\u00d752=|34\u00d730=|35\u00d753=|80\u00d784=\n77\u00d752=|49\u00d754=|37\u00d762=|75\u00d758=|14\u00d771=\n84\u00d728=|20\u00d716=|73\u00d745=|74\u00d788=|34\u00d711=\n96\u00d745=|36\u00d775=|28\u00d736=|40\u00d749=|57\u00d793=\n62\u00d718=|44\u00d743=|70\u00d785=|73\u00d793=|89\u00d756=\n76\u00d721=|31\u00d772=|14\u00d763=|33\u00d717=|41\u00d710=\n68\u00d771=|74\u00d789=|10\u00d733=|74\u00d784=|94\u00d736=\n10\u00d782=|18\u00d767=|43\u00d754=|45\u00d723=|77\u00d782=\n\"@ -split \"`r?`n\" | Where-Object { $_.Length -gt 0 }\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nif ($table.Rows.Count -ne $newRows.Count) {\n    throw \"Table has $($table.Rows.Count) rows, expected $($newRows.Count).\"\n}\n\n$updated = 0\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $rowValues = $newRows[$r - 1] -split \"\\|\"\n    if ($table.Columns.Count -ne $rowValues.Count) {\n        throw \"Table has $($table.Columns.Count) columns, expected $($rowValues.Count).\"\n    }\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        # Cell.Range.Text replaces only the run text; the existing run/\n        # paragraph formatting (font, size, alignment) is preserved.\n        $table.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n        $updated++\n    }\n}\n\nWrite-Output \"Updated $updated cells\"\n"}
